$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2").Value = "74.828.56"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3: 'Ethereum'
$ws.Range("D3").Value = "2.814.46"
$ws.Range("E3").Value = "  +7.57%  "

# Row 4: 'TetherUSD'
$ws.Range("E4").Value = "  +0.16%  "

# Row 5: 'Solana'
$ws.Range("D5").Value = "188.05"
$ws.Range("E5").Value = "  +0.24%  "

# Row 6: 'BNB'
$ws.Range("D6").Value = "596.99"
$ws.Range("E6").Value = "  +1.44%  "

# Row 7: 'USDC'
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8: 'XRP'
$ws.Range("D8").Value = "0.557"
$ws.Range("E8").Value = "  +3.27%  "

# Row 9: 'Dogecoin'
$ws.Range("E9").Value = "  -7.02%  "

# Row 10: 'LidoStakedEther'
$ws.Range("D10").Value = "2.812.31"
$ws.Range("E10").Value = "  +8.11%  "

# Row 11: 'TRON'
$ws.Range("E11").Value = "  -1.06%  "

# Row 12: 'Cardano'
$ws.Range("E12").Value = "  +1.50%  "

# Row 13: 'Toncoin'
$ws.Range("E13").Value = "  +0.92%  "

# Row 14: 'WrappedliquidstakedEther2.0'
$ws.Range("D14").Value = "3.327.11"
$ws.Range("E14").Value = "  +7.89%  "

# Row 15: 'WrappedBTC'
$ws.Range("D15").Value = "74.843.19"
$ws.Range("E15").Value = "  +0.61%  "

# Row 16: 'Avalanche'
$ws.Range("D16").Value = "27.01"
$ws.Range("E16").Value = "  +2.13%  "

# Row 17: 'ShibaInu'
$ws.Range("E17").Value = "  -3.37%  "

# Row 18: 'WrappedEther'
$ws.Range("D18").Value = "2.806.01"
$ws.Range("E18").Value = "  +6.50%  "

# Row 19: 'Uniswap'
$ws.Range("E19").Value = "  -2.18%  "

# Row 20: 'Chainlink'
$ws.Range("D20").Value = "12.32"
$ws.Range("E20").Value = "  +3.67%  "

# Row 21: 'BitcoinCash'
$ws.Range("D21").Value = "374.48"
$ws.Range("E21").Value = "  -0.38%  "

# Row 22: 'SuiNetwork'
$ws.Range("D22").Value = "2.25"
$ws.Range("E22").Value = "  -1.84%  "

# Row 23: 'Polkadot'
$ws.Range("E23").Value = "  -0.30%  "

# Row 24: 'LEO'
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "6.16"
$ws.Range("E24").Value = "  -1.03%  "

# Row 25: 'Dai'
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.00%  "

# Row 26: 'Litecoin'
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "70.59"
$ws.Range("E26").Value = "  +0.48%  "

# Row 27: 'WrappedeETH'
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.956.92"
$ws.Range("E27").Value = "  +8.20%  "

# Row 28: 'NEARProtocol'
$ws.Range("B28").Value = "NEARProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D28").Value = "4.15"
$ws.Range("E28").Value = "  -1.05%  "

# Row 29: 'Aptos'
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "9.56"
$ws.Range("E29").Value = "  +1.90%  "

# Row 30: 'PEPE'
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0000102"
$ws.Range("E30").Value = "  +7.34%  "

# Row 31: 'Binance-PegBSC-USD'
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.12%  "

# Row 32: 'Bittensor'
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "514.72"
$ws.Range("E32").Value = "  +0.59%  "

# Row 33: 'Fetch.AI'
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.38"
$ws.Range("E33").Value = "  -2.06%  "

# Row 34: 'InternetComputer(DFINITY)'
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "7.85"
$ws.Range("E34").Value = "  -1.97%  "

# Row 35: 'PancakeSwap'
$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").Value = "1.79"
$ws.Range("E35").Value = "  +1.68%  "

# Row 36: 'FirstDigitalUSD'
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.16%  "

# Row 37: 'Monero'
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "163.05"
$ws.Range("E37").Value = "  +1.98%  "

# Row 38: 'EthereumClassic'
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "20.08"
$ws.Range("E38").Value = "  +4.33%  "

# Row 39: 'Kaspa'
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.119"
$ws.Range("E39").Value = "  -3.68%  "

# Row 40: 'WhiteBITCoin'
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "19.31"
$ws.Range("E40").Value = "  -0.17%  "

# Row 41: 'Aave'
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "182.82"
$ws.Range("E41").Value = "  +15.63%  "

# Row 42: 'USDe'
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.04%  "

# Row 43: 'RenderToken'
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "5.03"
$ws.Range("E43").Value = "  +1.39%  "

# Row 44: 'PolygonEcosystemToken'
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "0.338"
$ws.Range("E44").Value = "  +3.22%  "

# Row 45: 'Stacks'
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "1.68"
$ws.Range("E45").Value = "  -1.24%  "

# Row 46: 'ImmutableX'
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "1.22"
$ws.Range("E46").Value = "  +2.62%  "

# Row 47: 'OKB'
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "39.87"
$ws.Range("E47").Value = "  +2.87%  "

# Row 48: 'dogwifhat'
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  -2.93%  "

# Row 49: 'Cronos'
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0858"
$ws.Range("E49").Value = "  -8.95%  "

# Row 50: 'ARBITRUM'
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "0.566"
$ws.Range("E50").Value = "  +7.39%  "

# Row 51: 'Filecoin'
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "3.73"
$ws.Range("E51").Value = "  +2.42%  "

Write-Host "Applied changes"